$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Update item names (column A) for rows 2-28 on both sheets
$ws1.Range("A2").Value = "MCER017"
$ws2.Range("A2").Value = "MCER017"
$ws1.Range("A3").Value = "MCER018"
$ws2.Range("A3").Value = "MCER018"
$ws1.Range("A4").Value = "MCER020"
$ws2.Range("A4").Value = "MCER020"
$ws1.Range("A5").Value = "MCER021"
$ws2.Range("A5").Value = "MCER021"
$ws1.Range("A6").Value = "MCER022"
$ws2.Range("A6").Value = "MCER022"
$ws1.Range("A7").Value = "MCER026"
$ws2.Range("A7").Value = "MCER026"
$ws1.Range("A8").Value = "MCER027"
$ws2.Range("A8").Value = "MCER027"
$ws1.Range("A9").Value = "MCER028"
$ws2.Range("A9").Value = "MCER028"
$ws1.Range("A10").Value = "MCER029"
$ws2.Range("A10").Value = "MCER029"
$ws1.Range("A11").Value = "MCER030"
$ws2.Range("A11").Value = "MCER030"
$ws1.Range("A12").Value = "MCER031"
$ws2.Range("A12").Value = "MCER031"
$ws1.Range("A13").Value = "MCER032"
$ws2.Range("A13").Value = "MCER032"
$ws1.Range("A14").Value = "MCER033"
$ws2.Range("A14").Value = "MCER033"
$ws1.Range("A15").Value = "MCER034"
$ws2.Range("A15").Value = "MCER034"
$ws1.Range("A16").Value = "MCER035"
$ws2.Range("A16").Value = "MCER035"
$ws1.Range("A17").Value = "MCER036"
$ws2.Range("A17").Value = "MCER036"
$ws1.Range("A18").Value = "MCER037"
$ws2.Range("A18").Value = "MCER037"
$ws1.Range("A19").Value = "MCER038"
$ws2.Range("A19").Value = "MCER038"
$ws1.Range("A20").Value = "MCER039"
$ws2.Range("A20").Value = "MCER039"
$ws1.Range("A21").Value = "MCER040"
$ws2.Range("A21").Value = "MCER040"
$ws1.Range("A22").Value = "MCER041"
$ws2.Range("A22").Value = "MCER041"
$ws1.Range("A23").Value = "MCER043"
$ws2.Range("A23").Value = "MCER043"
$ws1.Range("A24").Value = "MCER046"
$ws2.Range("A24").Value = "MCER046"
$ws1.Range("A25").Value = "MCER047"
$ws2.Range("A25").Value = "MCER047"
$ws1.Range("A26").Value = "MCER051"
$ws2.Range("A26").Value = "MCER051"
$ws1.Range("A27").Value = "MCER057"
$ws2.Range("A27").Value = "MCER057"
$ws1.Range("A28").Value = "MCER067"
$ws2.Range("A28").Value = "MCER067"

# Update minimum_stock (B) and maximum_stock (C) values for rows 2-28
$ws1.Range("B2").Value = 92
$ws1.Range("C2").Value = 793
$ws2.Range("B2").Value = 7
$ws2.Range("C2").Value = 642
$ws1.Range("B3").Value = 26
$ws1.Range("C3").Value = 836
$ws2.Range("B3").Value = 66
$ws2.Range("C3").Value = 288
$ws1.Range("B4").Value = 68
$ws1.Range("C4").Value = 866
$ws2.Range("B4").Value = 38
$ws2.Range("C4").Value = 373
$ws1.Range("B5").Value = 43
$ws1.Range("C5").Value = 555
$ws2.Range("B5").Value = 54
$ws2.Range("C5").Value = 546
$ws1.Range("B6").Value = 78
$ws1.Range("C6").Value = 418
$ws2.Range("B6").Value = 89
$ws2.Range("C6").Value = 653
$ws1.Range("B7").Value = 32
$ws1.Range("C7").Value = 292
$ws2.Range("B7").Value = 49
$ws2.Range("C7").Value = 707
$ws1.Range("B8").Value = 83
$ws1.Range("C8").Value = 636
$ws2.Range("B8").Value = 36
$ws2.Range("C8").Value = 653
$ws1.Range("B9").Value = 92
$ws1.Range("C9").Value = 765
$ws2.Range("B9").Value = 38
$ws2.Range("C9").Value = 610
$ws1.Range("B10").Value = 36
$ws1.Range("C10").Value = 239
$ws2.Range("B10").Value = 48
$ws2.Range("C10").Value = 828
$ws1.Range("B11").Value = 77
$ws1.Range("C11").Value = 730
$ws2.Range("B11").Value = 54
$ws2.Range("C11").Value = 767
$ws1.Range("B12").Value = 30
$ws1.Range("C12").Value = 440
$ws2.Range("B12").Value = 22
$ws2.Range("C12").Value = 905
$ws1.Range("B13").Value = 33
$ws1.Range("C13").Value = 621
$ws2.Range("B13").Value = 88
$ws2.Range("C13").Value = 282
$ws1.Range("B14").Value = 53
$ws1.Range("C14").Value = 363
$ws2.Range("B14").Value = 8
$ws2.Range("C14").Value = 414
$ws1.Range("B15").Value = 48
$ws1.Range("C15").Value = 667
$ws2.Range("B15").Value = 66
$ws2.Range("C15").Value = 633
$ws1.Range("B16").Value = 71
$ws1.Range("C16").Value = 652
$ws2.Range("B16").Value = 25
$ws2.Range("C16").Value = 613
$ws1.Range("B17").Value = 25
$ws1.Range("C17").Value = 751
$ws2.Range("B17").Value = 16
$ws2.Range("C17").Value = 441
$ws1.Range("B18").Value = 62
$ws1.Range("C18").Value = 706
$ws2.Range("B18").Value = 78
$ws2.Range("C18").Value = 577
$ws1.Range("B19").Value = 87
$ws1.Range("C19").Value = 276
$ws2.Range("B19").Value = 9
$ws2.Range("C19").Value = 623
$ws1.Range("B20").Value = 90
$ws1.Range("C20").Value = 159
$ws2.Range("B20").Value = 42
$ws2.Range("C20").Value = 207
$ws1.Range("B21").Value = 44
$ws1.Range("C21").Value = 773
$ws2.Range("B21").Value = 86
$ws2.Range("C21").Value = 243
$ws1.Range("B22").Value = 24
$ws1.Range("C22").Value = 182
$ws2.Range("B22").Value = 2
$ws2.Range("C22").Value = 796
$ws1.Range("B23").Value = 97
$ws1.Range("C23").Value = 528
$ws2.Range("B23").Value = 37
$ws2.Range("C23").Value = 747
$ws1.Range("B24").Value = 77
$ws1.Range("C24").Value = 717
$ws2.Range("B24").Value = 71
$ws2.Range("C24").Value = 161
$ws1.Range("B25").Value = 79
$ws1.Range("C25").Value = 372
$ws2.Range("B25").Value = 88
$ws2.Range("C25").Value = 258
$ws1.Range("B26").Value = 52
$ws1.Range("C26").Value = 126
$ws2.Range("B26").Value = 44
$ws2.Range("C26").Value = 993
$ws1.Range("B27").Value = 95
$ws1.Range("C27").Value = 136
$ws2.Range("B27").Value = 66
$ws2.Range("C27").Value = 204
$ws1.Range("B28").Value = 74
$ws1.Range("C28").Value = 238
$ws2.Range("B28").Value = 100
$ws2.Range("C28").Value = 636

# Remove now-unused rows 29-31 (previously held item_29..item_32 style rows)
$ws1.Rows("29:31").Delete()
$ws2.Rows("29:31").Delete()
